$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateAccount")

# Row 7 - plain text row (no hyperlink)
$ws.Range("A7").Value = "carla"
$ws.Range("B7").Value = "brook"
$ws.Range("C7").Value = "brook567@ymail.com"
$ws.Range("D7").Value = "9876fdsa"

# Row 8 - Jason/Thomas, with a mailto hyperlink on C8
$ws.Range("A8").Value = "Jason"
$ws.Range("B8").Value = "Thomas"
$ws.Range("D8").Value = "1334348jdfhd"

$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:mike34@yahoo.com", "", "", "jas@")
$ws.Range("C8").Style = "Hyperlink"

$ws.Range("B8").Select()
